$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at D and E (shifting batsman..sr from D..I to F..K)
$ws.Range("D1:E1").EntireColumn.Insert()

# New header values
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"

# New data values
$ws.Range("D2").Value = "Delhi Capitals"
$ws.Range("E2").Value = "Mumbai Indians"
